# Update time_taken (column F) values on the "data" sheet to reflect the
# new panel query run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$timeTaken = @(
    "2021-10-05 14:20:35.425339",
    "2021-10-05 14:20:35.425347",
    "2021-10-05 14:20:35.425351",
    "2021-10-05 14:20:35.425353",
    "2021-10-05 14:20:35.425356",
    "2021-10-05 14:20:35.425359",
    "2021-10-05 14:20:35.425361",
    "2021-10-05 14:20:35.425364",
    "2021-10-05 14:20:35.425367",
    "2021-10-05 14:20:35.425369",
    "2021-10-05 14:20:35.425372",
    "2021-10-05 14:20:35.425374",
    "2021-10-05 14:20:35.425377",
    "2021-10-05 14:20:35.425379",
    "2021-10-05 14:20:35.425382",
    "2021-10-05 14:20:35.425384",
    "2021-10-05 14:20:35.425387",
    "2021-10-05 14:20:35.425390",
    "2021-10-05 14:20:35.425392",
    "2021-10-05 14:20:35.425395",
    "2021-10-05 14:20:35.425397",
    "2021-10-05 14:20:35.425400",
    "2021-10-05 14:20:35.425402",
    "2021-10-05 14:20:35.425405",
    "2021-10-05 14:20:35.425407",
    "2021-10-05 14:20:35.425410",
    "2021-10-05 14:20:35.425412",
    "2021-10-05 14:20:35.425415",
    "2021-10-05 14:20:35.425417",
    "2021-10-05 14:20:35.425420",
    "2021-10-05 14:20:35.425422",
    "2021-10-05 14:20:35.425425",
    "2021-10-05 14:20:35.425428",
    "2021-10-05 14:20:35.425430",
    "2021-10-05 14:20:35.425433",
    "2021-10-05 14:20:35.425435",
    "2021-10-05 14:20:35.425438",
    "2021-10-05 14:20:35.425440",
    "2021-10-05 14:20:35.425443",
    "2021-10-05 14:20:35.425445",
    "2021-10-05 14:20:35.425448",
    "2021-10-05 14:20:35.425451",
    "2021-10-05 14:20:35.425453",
    "2021-10-05 14:20:35.425456",
    "2021-10-05 14:20:35.425458",
    "2021-10-05 14:20:35.425461",
    "2021-10-05 14:20:35.425463",
    "2021-10-05 14:20:35.425466",
    "2021-10-05 14:20:35.425468",
    "2021-10-05 14:20:35.425471",
    "2021-10-05 14:20:35.425473",
    "2021-10-05 14:20:35.425476",
    "2021-10-05 14:20:35.425478",
    "2021-10-05 14:20:35.425481",
    "2021-10-05 14:20:35.425483",
    "2021-10-05 14:20:35.425486",
    "2021-10-05 14:20:35.425488",
    "2021-10-05 14:20:35.425491",
    "2021-10-05 14:20:35.425493",
    "2021-10-05 14:20:35.425496",
    "2021-10-05 14:20:35.425499",
    "2021-10-05 14:20:35.425501",
    "2021-10-05 14:20:35.425504",
    "2021-10-05 14:20:35.425506",
    "2021-10-05 14:20:35.425510",
    "2021-10-05 14:20:35.425513",
    "2021-10-05 14:20:35.425515",
    "2021-10-05 14:20:35.425518",
    "2021-10-05 14:20:35.425520",
    "2021-10-05 14:20:35.425523",
    "2021-10-05 14:20:35.425525",
    "2021-10-05 14:20:35.425528",
    "2021-10-05 14:20:35.425530",
    "2021-10-05 14:20:35.425533",
    "2021-10-05 14:20:35.425535",
    "2021-10-05 14:20:35.425538",
    "2021-10-05 14:20:35.425542",
    "2021-10-05 14:20:35.425545",
    "2021-10-05 14:20:35.425548",
    "2021-10-05 14:20:35.425551",
    "2021-10-05 14:20:35.425553",
    "2021-10-05 14:20:35.425555",
    "2021-10-05 14:20:35.425558",
    "2021-10-05 14:20:35.425560",
    "2021-10-05 14:20:35.425563",
    "2021-10-05 14:20:35.425565",
    "2021-10-05 14:20:35.425568",
    "2021-10-05 14:20:35.425570",
    "2021-10-05 14:20:35.425573",
    "2021-10-05 14:20:35.425575",
    "2021-10-05 14:20:35.425578",
    "2021-10-05 14:20:35.425580",
    "2021-10-05 14:20:35.425584",
    "2021-10-05 14:20:35.425587",
    "2021-10-05 14:20:35.425589",
    "2021-10-05 14:20:35.425592",
    "2021-10-05 14:20:35.425594",
    "2021-10-05 14:20:35.425597",
    "2021-10-05 14:20:35.425599",
    "2021-10-05 14:20:35.425602",
    "2021-10-05 14:20:35.425604",
    "2021-10-05 14:20:35.425607",
    "2021-10-05 14:20:35.425609",
    "2021-10-05 14:20:35.425612",
    "2021-10-05 14:20:35.425614",
    "2021-10-05 14:20:35.425617",
    "2021-10-05 14:20:35.425619",
    "2021-10-05 14:20:35.425622",
    "2021-10-05 14:20:35.425626",
    "2021-10-05 14:20:35.425629",
    "2021-10-05 14:20:35.425632",
    "2021-10-05 14:20:35.425634",
    "2021-10-05 14:20:35.425637",
    "2021-10-05 14:20:35.425640",
    "2021-10-05 14:20:35.425642",
    "2021-10-05 14:20:35.425645",
    "2021-10-05 14:20:35.425647",
    "2021-10-05 14:20:35.425650",
    "2021-10-05 14:20:35.425652",
    "2021-10-05 14:20:35.425655",
    "2021-10-05 14:20:35.425658",
    "2021-10-05 14:20:35.425661",
    "2021-10-05 14:20:35.425663",
    "2021-10-05 14:20:35.425666",
    "2021-10-05 14:20:35.425668",
    "2021-10-05 14:20:35.425671",
    "2021-10-05 14:20:35.425673",
    "2021-10-05 14:20:35.425676",
    "2021-10-05 14:20:35.425680",
    "2021-10-05 14:20:35.425684",
    "2021-10-05 14:20:35.425686",
    "2021-10-05 14:20:35.425689",
    "2021-10-05 14:20:35.425691",
    "2021-10-05 14:20:35.425694",
    "2021-10-05 14:20:35.425697",
    "2021-10-05 14:20:35.425699",
    "2021-10-05 14:20:35.425702",
    "2021-10-05 14:20:35.425704",
    "2021-10-05 14:20:35.425707",
    "2021-10-05 14:20:35.425710",
    "2021-10-05 14:20:35.425712",
    "2021-10-05 14:20:35.425715",
    "2021-10-05 14:20:35.425717",
    "2021-10-05 14:20:35.425720",
    "2021-10-05 14:20:35.425723",
    "2021-10-05 14:20:35.425725",
    "2021-10-05 14:20:35.425728",
    "2021-10-05 14:20:35.425731",
    "2021-10-05 14:20:35.425734",
    "2021-10-05 14:20:35.425736",
    "2021-10-05 14:20:35.425739"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

# Add a new "metadata" sheet after the existing "data" sheet, describing the
# panel query that produced this workbook.
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Reuse the bold/bordered header style already used on the "data" sheet for
# the metadata header row and the A2 index cell.
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$meta.Range("A2").PasteSpecial(-4122)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Growth failure in early childhood"
$meta.Range("C2").Value = 473
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.82"
$meta.Range("E2").Value = "2021-10-01T14:36:27.027770Z"
$meta.Range("F2").Value = "2021-10-05 14:20:35.422041"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/473/?format=json"
